# Updates recalculated noise_threshold / First Rise Point results after
# adding the configurable zero_before_threshold parameter.
#
# Columns:
#   C = (noise threshold) index/count
#   E = recalculated fraction/value tied to the new threshold logic
#   G = recalculated rise-point index
#
# The same C/E updates apply to rows 2-6 on every "Step3_DataPts_*" sheet;
# only the G column values differ per sheet (per percentile).

$wb = $excel.ActiveWorkbook

# Row data shared by all four Step3_DataPts sheets for columns C and E.
$rows = @(2, 3, 4, 5, 6)
$newC = @{ 2 = 88; 3 = 88; 4 = 43; 5 = 88; 6 = 43 }
$newE = @{
    2 = 0.0005376466737016518
    3 = 0.01635448674696752
    4 = 0.03147651671560622
    5 = 0.02995609021605257
    6 = 0.01514829473258368
}

# Per-sheet new G column values.
$newG = @{
    "Step3_DataPts_0.5" = @{ 2 = 36; 3 = 42; 4 = 55; 5 = 43; 6 = 55 }
    "Step3_DataPts_0.7" = @{ 2 = 56; 3 = 56; 4 = 72; 5 = 57; 6 = 71 }
    "Step3_DataPts_0.8" = @{ 2 = 67; 3 = 67; 4 = 99; 5 = 67; 6 = 90 }
    "Step3_DataPts_0.9" = @{ 2 = 89; 3 = 89; 4 = 122; 5 = 90; 6 = 121 }
}

foreach ($sheetName in $newG.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $gValues = $newG[$sheetName]

    foreach ($r in $rows) {
        $ws.Range("C$r").Value = $newC[$r]
        $ws.Range("E$r").Value = $newE[$r]
        $ws.Range("G$r").Value = $gValues[$r]
    }
}
